$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Fix up the "保險" (insurance) sheet: add the missing field-name
#    header columns (company/name/owner/property_category/category/
#    date/legislator_name/legislator_id/source_file/index) and put
#    the correct per-row values (company name, policy name, owner)
#    in B/C/D while filling the newly added columns E-K, matching the
#    layout used by the other property sheets (e.g. "股票").
# ------------------------------------------------------------------
$insurance = $wb.Worksheets.Item("保險")

# --- Row 1 (header labels), columns B-D already exist with the bold
#     bordered header style (style index 1) - just fix their text.
$insurance.Cells.Item(1,2).Value = "company"
$insurance.Cells.Item(1,3).Value = "name"
$insurance.Cells.Item(1,4).Value = "owner"

# New header cells E1:K1 - copy formatting from the existing D1 header
# cell (bold, centered, bordered) then set their text.
$headerLabels = @("property_category","category","date","legislator_name","legislator_id","source_file","index")
for ($i = 0; $i -lt $headerLabels.Length; $i++) {
    $col = 5 + $i
    $insurance.Cells.Item(1,4).Copy()
    $insurance.Cells.Item(1,$col).PasteSpecial(-4122)
    $insurance.Cells.Item(1,$col).Value = $headerLabels[$i]
}

# --- Row 2 (富邦人壽 / 安泰人壽靈活理財變額保險甲型) ---
$insurance.Cells.Item(2,2).Value = "富邦人壽"
$insurance.Cells.Item(2,3).Value = "安泰人壽靈活理財變額保險甲型"
$insurance.Cells.Item(2,4).Value = "賴士葆"

$insurance.Cells.Item(2,4).Copy()
$insurance.Cells.Item(2,5).PasteSpecial(-4122)
$insurance.Cells.Item(2,5).Value = "insurance"

$insurance.Cells.Item(2,4).Copy()
$insurance.Cells.Item(2,6).PasteSpecial(-4122)
$insurance.Cells.Item(2,6).Value = "normal"

$insurance.Cells.Item(2,4).Copy()
$insurance.Cells.Item(2,7).PasteSpecial(-4122)
$insurance.Cells.Item(2,7).NumberFormat = "@"
$insurance.Cells.Item(2,7).Value = "2012-04-19"

$insurance.Cells.Item(2,4).Copy()
$insurance.Cells.Item(2,8).PasteSpecial(-4122)
$insurance.Cells.Item(2,8).Value = "賴士葆"

$insurance.Cells.Item(2,4).Copy()
$insurance.Cells.Item(2,9).PasteSpecial(-4122)
$insurance.Cells.Item(2,9).Value = 866

$insurance.Cells.Item(2,4).Copy()
$insurance.Cells.Item(2,10).PasteSpecial(-4122)
$insurance.Cells.Item(2,10).Value = "tmp9edb1"

$insurance.Cells.Item(2,4).Copy()
$insurance.Cells.Item(2,11).PasteSpecial(-4122)
$insurance.Cells.Item(2,11).Value = 102

# --- Row 3 (國華人壽 / 國華人壽終身壽險) ---
$insurance.Cells.Item(3,2).Value = "國華人壽"
$insurance.Cells.Item(3,3).Value = "國華人壽終身壽險"
$insurance.Cells.Item(3,4).Value = "賴士葆"

$insurance.Cells.Item(3,4).Copy()
$insurance.Cells.Item(3,5).PasteSpecial(-4122)
$insurance.Cells.Item(3,5).Value = "insurance"

$insurance.Cells.Item(3,4).Copy()
$insurance.Cells.Item(3,6).PasteSpecial(-4122)
$insurance.Cells.Item(3,6).Value = "normal"

$insurance.Cells.Item(3,4).Copy()
$insurance.Cells.Item(3,7).PasteSpecial(-4122)
$insurance.Cells.Item(3,7).NumberFormat = "@"
$insurance.Cells.Item(3,7).Value = "2012-04-19"

$insurance.Cells.Item(3,4).Copy()
$insurance.Cells.Item(3,8).PasteSpecial(-4122)
$insurance.Cells.Item(3,8).Value = "賴士葆"

$insurance.Cells.Item(3,4).Copy()
$insurance.Cells.Item(3,9).PasteSpecial(-4122)
$insurance.Cells.Item(3,9).Value = 866

$insurance.Cells.Item(3,4).Copy()
$insurance.Cells.Item(3,10).PasteSpecial(-4122)
$insurance.Cells.Item(3,10).Value = "tmp9edb1"

$insurance.Cells.Item(3,4).Copy()
$insurance.Cells.Item(3,11).PasteSpecial(-4122)
$insurance.Cells.Item(3,11).Value = 103

# ------------------------------------------------------------------
# 2) Remove the stray "債務" (debt) sheet entirely - its content was
#    a leftover mix of investment-section labels and is dropped.
# ------------------------------------------------------------------
$debt = $wb.Worksheets.Item("債務")
$debt.Delete()
